# Update the "Förändrad" (last-changed) date in column C for rows 2-9
# from 2023-11-03 (45233) to 2023-11-13 (45243), as produced by the
# automatic data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45243
}
